# This workbook is a weekly price log for "Brócoli" (Femacal de La Calera).
# Each week two rows (a "Primera" and a "Segunda" quality row) are inserted
# into the existing table. This edit inserts two such weekly pairs:
#   - a pair dated 2022-07-06 (serial 44748), inserted at row 642
#   - a pair dated 2022-07-05 (serial 44747), inserted a bit further down
#     (at what becomes row 698 after the first insertion)
# All the rows that were already in the table simply get pushed down to make
# room - none of their own data changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new weekly pair (Primera/Segunda) at row 642 ---------
$ws.Rows("642:643").Insert()

$ws.Cells.Item(642, 1).Value = 3
$ws.Cells.Item(642, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(642, 3).Value = "Coquimbo"
$ws.Cells.Item(642, 4).Value = 44748
$ws.Cells.Item(642, 5).Value = 5
$ws.Cells.Item(642, 6).Value = 100112023
$ws.Cells.Item(642, 7).Value = "Brócoli"
$ws.Cells.Item(642, 8).Value = "Sin especificar"
$ws.Cells.Item(642, 9).Value = "Primera"
$ws.Cells.Item(642, 10).Value = 1900
$ws.Cells.Item(642, 11).Value = 800
$ws.Cells.Item(642, 12).Value = 850
$ws.Cells.Item(642, 13).Value = 825
$ws.Cells.Item(642, 14).Value = "`$/unidad"
$ws.Cells.Item(642, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(642, 16).Value = 825
$ws.Cells.Item(642, 17).Value = 1
$ws.Cells.Item(642, 18).Value = "Hortaliza"

$ws.Cells.Item(643, 1).Value = 3
$ws.Cells.Item(643, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(643, 3).Value = "Coquimbo"
$ws.Cells.Item(643, 4).Value = 44748
$ws.Cells.Item(643, 5).Value = 5
$ws.Cells.Item(643, 6).Value = 100112023
$ws.Cells.Item(643, 7).Value = "Brócoli"
$ws.Cells.Item(643, 8).Value = "Sin especificar"
$ws.Cells.Item(643, 9).Value = "Segunda"
$ws.Cells.Item(643, 10).Value = 900
$ws.Cells.Item(643, 11).Value = 700
$ws.Cells.Item(643, 12).Value = 700
$ws.Cells.Item(643, 13).Value = 700
$ws.Cells.Item(643, 14).Value = "`$/unidad"
$ws.Cells.Item(643, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(643, 16).Value = 700
$ws.Cells.Item(643, 17).Value = 1
$ws.Cells.Item(643, 18).Value = "Hortaliza"

# --- Insert second new weekly pair (Primera/Segunda) at row 698 --------
$ws.Rows("698:699").Insert()

$ws.Cells.Item(698, 1).Value = 3
$ws.Cells.Item(698, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(698, 3).Value = "Coquimbo"
$ws.Cells.Item(698, 4).Value = 44747
$ws.Cells.Item(698, 5).Value = 5
$ws.Cells.Item(698, 6).Value = 100112023
$ws.Cells.Item(698, 7).Value = "Brócoli"
$ws.Cells.Item(698, 8).Value = "Sin especificar"
$ws.Cells.Item(698, 9).Value = "Primera"
$ws.Cells.Item(698, 10).Value = 2750
$ws.Cells.Item(698, 11).Value = 800
$ws.Cells.Item(698, 12).Value = 850
$ws.Cells.Item(698, 13).Value = 825
$ws.Cells.Item(698, 14).Value = "`$/unidad"
$ws.Cells.Item(698, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(698, 16).Value = 825
$ws.Cells.Item(698, 17).Value = 1
$ws.Cells.Item(698, 18).Value = "Hortaliza"

$ws.Cells.Item(699, 1).Value = 3
$ws.Cells.Item(699, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(699, 3).Value = "Coquimbo"
$ws.Cells.Item(699, 4).Value = 44747
$ws.Cells.Item(699, 5).Value = 5
$ws.Cells.Item(699, 6).Value = 100112023
$ws.Cells.Item(699, 7).Value = "Brócoli"
$ws.Cells.Item(699, 8).Value = "Sin especificar"
$ws.Cells.Item(699, 9).Value = "Segunda"
$ws.Cells.Item(699, 10).Value = 1700
$ws.Cells.Item(699, 11).Value = 700
$ws.Cells.Item(699, 12).Value = 700
$ws.Cells.Item(699, 13).Value = 700
$ws.Cells.Item(699, 14).Value = "`$/unidad"
$ws.Cells.Item(699, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(699, 16).Value = 700
$ws.Cells.Item(699, 17).Value = 1
$ws.Cells.Item(699, 18).Value = "Hortaliza"
